$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NextBus3")

# Remove the "Right_BusStopCode" column (M) - data synced from source no longer includes it
$ws.Range("M1").EntireColumn.Delete()

# Remove now-stale rows 11-15 (refreshed source only has 9 data rows)
$ws.Range("A11:A15").EntireRow.Delete()

# Refresh remaining data rows (2-10) with latest synced values
$ws.Cells.Item(2, 1).Value2 = 'NextBus3'
$ws.Cells.Item(2, 2).Value2 = 52
$ws.Cells.Item(2, 3).Value2 = 53009
$ws.Cells.Item(2, 4).Value2 = 'Bishan Int'
$ws.Cells.Item(2, 5).Value2 = 'SBST'
$ws.Cells.Item(2, 6).Value2 = 45685.00008101852
$ws.Cells.Item(2, 7).Value2 = 53009
$ws.Cells.Item(2, 8).Value2 = 'WAB'
$ws.Cells.Item(2, 9).Value2 = 'SEA'
$ws.Cells.Item(2, 10).Value2 = 0
$ws.Cells.Item(2, 11).Value2 = 28009
$ws.Cells.Item(2, 12).Value2 = 'DD'
$ws.Cells.Item(2, 13).Value2 = 612
$ws.Cells.Item(2, 14).Value2 = 15
$ws.Cells.Item(2, 15).Value2 = 613
$ws.Cells.Item(2, 16).Value2 = 16
$ws.Cells.Item(2, 17).Value2 = 619
$ws.Cells.Item(2, 18).Value2 = 15
$ws.Cells.Item(2, 19).Value2 = 12101
$ws.Cells.Item(2, 20).Value2 = 'Ngee Ann Poly'
$ws.Cells.Item(3, 1).Value2 = 'NextBus3'
$ws.Cells.Item(3, 2).Value2 = 184
$ws.Cells.Item(3, 3).Value2 = 44989
$ws.Cells.Item(3, 4).Value2 = 'Gali Batu Ter'
$ws.Cells.Item(3, 5).Value2 = 'SMRT'
$ws.Cells.Item(3, 6).Value2 = 45684.99182870371
$ws.Cells.Item(3, 7).Value2 = 44989
$ws.Cells.Item(3, 8).Value2 = 'WAB'
$ws.Cells.Item(3, 9).Value2 = 'SEA'
$ws.Cells.Item(3, 10).Value2 = 1
$ws.Cells.Item(3, 11).Value2 = 44989
$ws.Cells.Item(3, 12).Value2 = 'DD'
$ws.Cells.Item(3, 13).Value2 = 638
$ws.Cells.Item(3, 14).Value2 = 2347
$ws.Cells.Item(3, 15).Value2 = 634
$ws.Cells.Item(3, 16).Value2 = 2345
$ws.Cells.Item(3, 17).Value2 = 612
$ws.Cells.Item(3, 18).Value2 = 2350
$ws.Cells.Item(3, 19).Value2 = 12101
$ws.Cells.Item(3, 20).Value2 = 'Ngee Ann Poly'
$ws.Cells.Item(4, 1).Value2 = 'NextBus3'
$ws.Cells.Item(4, 2).Value2 = 74
$ws.Cells.Item(4, 3).Value2 = 11379
$ws.Cells.Item(4, 4).Value2 = 'Buona Vista Ter'
$ws.Cells.Item(4, 5).Value2 = 'SBST'
$ws.Cells.Item(4, 6).Value2 = 45684.98525462963
$ws.Cells.Item(4, 7).Value2 = 11379
$ws.Cells.Item(4, 8).Value2 = 'WAB'
$ws.Cells.Item(4, 9).Value2 = 'SEA'
$ws.Cells.Item(4, 10).Value2 = 1
$ws.Cells.Item(4, 11).Value2 = 64009
$ws.Cells.Item(4, 12).Value2 = 'DD'
$ws.Cells.Item(4, 13).Value2 = 605
$ws.Cells.Item(4, 14).Value2 = 30
$ws.Cells.Item(4, 15).Value2 = 559
$ws.Cells.Item(4, 16).Value2 = 24
$ws.Cells.Item(4, 17).Value2 = 609
$ws.Cells.Item(4, 18).Value2 = 27
$ws.Cells.Item(4, 19).Value2 = 12109
$ws.Cells.Item(4, 20).Value2 = 'Opp Ngee Ann Poly'
$ws.Cells.Item(5, 1).Value2 = 'NextBus3'
$ws.Cells.Item(5, 2).Value2 = 154
$ws.Cells.Item(5, 3).Value2 = 82009
$ws.Cells.Item(5, 4).Value2 = 'Eunos Int'
$ws.Cells.Item(5, 5).Value2 = 'SBST'
$ws.Cells.Item(5, 6).Value2 = 45684.98541666667
$ws.Cells.Item(5, 7).Value2 = 82009
$ws.Cells.Item(5, 8).Value2 = 'WAB'
$ws.Cells.Item(5, 9).Value2 = 'SEA'
$ws.Cells.Item(5, 10).Value2 = 1
$ws.Cells.Item(5, 11).Value2 = 22009
$ws.Cells.Item(5, 12).Value2 = 'SD'
$ws.Cells.Item(5, 13).Value2 = 601
$ws.Cells.Item(5, 14).Value2 = 57
$ws.Cells.Item(5, 15).Value2 = 637
$ws.Cells.Item(5, 16).Value2 = 55
$ws.Cells.Item(5, 17).Value2 = 604
$ws.Cells.Item(5, 18).Value2 = 53
$ws.Cells.Item(5, 19).Value2 = 12101
$ws.Cells.Item(5, 20).Value2 = 'Ngee Ann Poly'
$ws.Cells.Item(6, 1).Value2 = 'NextBus3'
$ws.Cells.Item(6, 2).Value2 = 154
$ws.Cells.Item(6, 3).Value2 = 22009
$ws.Cells.Item(6, 4).Value2 = 'Boon Lay Int'
$ws.Cells.Item(6, 5).Value2 = 'SBST'
$ws.Cells.Item(6, 6).Value2 = 45684.99453703704
$ws.Cells.Item(6, 7).Value2 = 22009
$ws.Cells.Item(6, 8).Value2 = 'WAB'
$ws.Cells.Item(6, 9).Value2 = 'SEA'
$ws.Cells.Item(6, 10).Value2 = 1
$ws.Cells.Item(6, 11).Value2 = 82009
$ws.Cells.Item(6, 12).Value2 = 'SD'
$ws.Cells.Item(6, 13).Value2 = 546
$ws.Cells.Item(6, 14).Value2 = 16
$ws.Cells.Item(6, 15).Value2 = 616
$ws.Cells.Item(6, 16).Value2 = 15
$ws.Cells.Item(6, 17).Value2 = 547
$ws.Cells.Item(6, 18).Value2 = 15
$ws.Cells.Item(6, 19).Value2 = 12109
$ws.Cells.Item(6, 20).Value2 = 'Opp Ngee Ann Poly'
$ws.Cells.Item(7, 1).Value2 = 'NextBus3'
$ws.Cells.Item(7, 2).Value2 = 61
$ws.Cells.Item(7, 3).Value2 = 43009
$ws.Cells.Item(7, 4).Value2 = 'Bt Batok Int'
$ws.Cells.Item(7, 5).Value2 = 'SMRT'
$ws.Cells.Item(7, 6).Value2 = 45684.98640046296
$ws.Cells.Item(7, 7).Value2 = 43009
$ws.Cells.Item(7, 8).Value2 = 'WAB'
$ws.Cells.Item(7, 9).Value2 = 'SEA'
$ws.Cells.Item(7, 10).Value2 = 1
$ws.Cells.Item(7, 11).Value2 = 82009
$ws.Cells.Item(7, 12).Value2 = 'DD'
$ws.Cells.Item(7, 13).Value2 = 645
$ws.Cells.Item(7, 14).Value2 = 108
$ws.Cells.Item(7, 15).Value2 = 652
$ws.Cells.Item(7, 16).Value2 = 110
$ws.Cells.Item(7, 17).Value2 = 642
$ws.Cells.Item(7, 18).Value2 = 109
$ws.Cells.Item(7, 19).Value2 = 12101
$ws.Cells.Item(7, 20).Value2 = 'Ngee Ann Poly'
$ws.Cells.Item(8, 1).Value2 = 'NextBus3'
$ws.Cells.Item(8, 2).Value2 = 52
$ws.Cells.Item(8, 3).Value2 = 28009
$ws.Cells.Item(8, 4).Value2 = 'Jurong East Int'
$ws.Cells.Item(8, 5).Value2 = 'SBST'
$ws.Cells.Item(8, 6).Value2 = 45684.99422453704
$ws.Cells.Item(8, 7).Value2 = 28009
$ws.Cells.Item(8, 8).Value2 = 'WAB'
$ws.Cells.Item(8, 9).Value2 = 'SEA'
$ws.Cells.Item(8, 10).Value2 = 1
$ws.Cells.Item(8, 11).Value2 = 53009
$ws.Cells.Item(8, 12).Value2 = 'SD'
$ws.Cells.Item(8, 13).Value2 = 623
$ws.Cells.Item(8, 14).Value2 = 23
$ws.Cells.Item(8, 15).Value2 = 625
$ws.Cells.Item(8, 16).Value2 = 21
$ws.Cells.Item(8, 17).Value2 = 627
$ws.Cells.Item(8, 18).Value2 = 22
$ws.Cells.Item(8, 19).Value2 = 12109
$ws.Cells.Item(8, 20).Value2 = 'Opp Ngee Ann Poly'
$ws.Cells.Item(9, 1).Value2 = 'NextBus3'
$ws.Cells.Item(9, 2).Value2 = 151
$ws.Cells.Item(9, 3).Value2 = 16009
$ws.Cells.Item(9, 4).Value2 = 'Kent Ridge Ter'
$ws.Cells.Item(9, 5).Value2 = 'SBST'
$ws.Cells.Item(9, 6).Value2 = 45684.99447916666
$ws.Cells.Item(9, 7).Value2 = 16009
$ws.Cells.Item(9, 8).Value2 = 'WAB'
$ws.Cells.Item(9, 9).Value2 = 'SEA'
$ws.Cells.Item(9, 10).Value2 = 1
$ws.Cells.Item(9, 11).Value2 = 64009
$ws.Cells.Item(9, 12).Value2 = 'SD'
$ws.Cells.Item(9, 13).Value2 = 635
$ws.Cells.Item(9, 14).Value2 = 2347
$ws.Cells.Item(9, 15).Value2 = 634
$ws.Cells.Item(9, 16).Value2 = 2351
$ws.Cells.Item(9, 17).Value2 = 639
$ws.Cells.Item(9, 18).Value2 = 2354
$ws.Cells.Item(9, 19).Value2 = 12109
$ws.Cells.Item(9, 20).Value2 = 'Opp Ngee Ann Poly'
$ws.Cells.Item(10, 1).Value2 = 'NextBus3'
$ws.Cells.Item(10, 2).Value2 = 75
$ws.Cells.Item(10, 3).Value2 = 10009
$ws.Cells.Item(10, 4).Value2 = 'Bt Merah Int'
$ws.Cells.Item(10, 5).Value2 = 'SMRT'
$ws.Cells.Item(10, 6).Value2 = 45684.99631944444
$ws.Cells.Item(10, 7).Value2 = 10009
$ws.Cells.Item(10, 8).Value2 = 'WAB'
$ws.Cells.Item(10, 9).Value2 = 'SEA'
$ws.Cells.Item(10, 10).Value2 = 0
$ws.Cells.Item(10, 11).Value2 = 44989
$ws.Cells.Item(10, 12).Value2 = 'SD'
$ws.Cells.Item(10, 13).Value2 = 548
$ws.Cells.Item(10, 14).Value2 = 2350
$ws.Cells.Item(10, 15).Value2 = 546
$ws.Cells.Item(10, 16).Value2 = 2350
$ws.Cells.Item(10, 17).Value2 = 552
$ws.Cells.Item(10, 18).Value2 = 2351
$ws.Cells.Item(10, 19).Value2 = 12109
$ws.Cells.Item(10, 20).Value2 = 'Opp Ngee Ann Poly'
